$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - column F values (想去人数)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 570
$ws1.Range("F3").Value = 188
$ws1.Range("F4").Value = 365
$ws1.Range("F5").Value = 414
$ws1.Range("F6").Value = 263
$ws1.Range("F7").Value = 2406
$ws1.Range("F8").Value = 414
$ws1.Range("F9").Value = 6269
$ws1.Range("F10").Value = 163
$ws1.Range("F11").Value = 403
$ws1.Range("F12").Value = 21

# Sheet "演出" (Performances) - column F values (想去人数)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 12
$ws2.Range("F3").Value = 16

# Sheet "全部类型" (All types) - column F values (想去人数)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 570
$ws4.Range("F3").Value = 188
$ws4.Range("F4").Value = 365
$ws4.Range("F5").Value = 414
$ws4.Range("F6").Value = 263
$ws4.Range("F7").Value = 12
$ws4.Range("F8").Value = 16
$ws4.Range("F9").Value = 2406
$ws4.Range("F10").Value = 414
$ws4.Range("F11").Value = 6269
$ws4.Range("F12").Value = 163
$ws4.Range("F13").Value = 403
$ws4.Range("F15").Value = 21

$wb.Save()
